# Applies the "added generate functionality to referral form and
# clientmonitoring" edit to the Client Monitoring Form document.

$d = $word.ActiveDocument

# 1) Campus heading placeholder -> example Campus
$d.Content.Find.Execute("try Campus", $true, $false, $false, $false, $false, `
    $true, 1, $false, "example Campus", 2) | Out-Null

# 2) Name field: drop the trailing "Ramada" (keep leading spaces + first two names)
$d.Content.Find.Execute("    John Vincent Ramada", $true, $false, $false, $false, $false, `
    $true, 1, $false, "    John Vincent ", 2) | Out-Null

# 3) Adviser field: "Lyster John" -> "John Vincent"
$d.Content.Find.Execute("     Lyster John", $true, $false, $false, $false, $false, `
    $true, 1, $false, "     John Vincent", 2) | Out-Null

# 4) Client monitoring table: update the generated-sample row (row 2) and
#    clear the second sample row (row 3), which becomes a blank template row.
$table = $d.Tables.Item(1)

$table.Cell(2, 1).Range.Text = "date"
$table.Cell(2, 2).Range.Text = "concern"
$table.Cell(2, 3).Range.Text = "signatures/98yGuuwbWatcEoCy2FQrK02DuJMpcbIoSol9snm9.jpg"
$table.Cell(2, 4).Range.Text = "action taken"
$table.Cell(2, 5).Range.Text = "recommendation"

$table.Cell(3, 1).Range.Text = ""
$table.Cell(3, 2).Range.Text = ""
$table.Cell(3, 3).Range.Text = ""
$table.Cell(3, 4).Range.Text = ""
$table.Cell(3, 5).Range.Text = ""
